$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 469 (shifts the existing rows 469-501 down to 470-502,
# preserving all of their original data and formatting).
$ws.Rows.Item(469).Insert()

# Populate the newly inserted row 469 with this week's new price record.
$ws.Range("A469").Value = 10
$ws.Range("B469").Value = "Vega Modelo de Temuco"
$ws.Range("C469").Value = "La Araucanía"
$ws.Range("D469").Value = 44516
$ws.Range("E469").Value = 9
$ws.Range("F469").Value = "Fruta"
$ws.Range("G469").Value = 100104
$ws.Range("H469").Value = "Frutos de pepita"
$ws.Range("I469").Value = 100104005
$ws.Range("J469").Value = "Pera"
$ws.Range("K469").Value = "Packham's Triumph"
$ws.Range("L469").Value = "Primera"
$ws.Range("M469").Value = 95
$ws.Range("N469").Value = 16000
$ws.Range("O469").Value = 16000
$ws.Range("P469").Value = 16000
$ws.Range("Q469").Value = "$/bandeja 18 kilos granel"
$ws.Range("R469").Value = "Región de O'Higgins"
$ws.Range("S469").Value = 889
$ws.Range("T469").Value = 18
